# Update "想去人数" (want-to-go count) figures in column F across the
# 展览 (sheet1), 演出 (sheet2) and 全部类型 (sheet4) tabs.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws1.Range("F2").Value = 143
$ws1.Range("F3").Value = 457
$ws1.Range("F6").Value = 11
$ws1.Range("F7").Value = 27
$ws1.Range("F9").Value = 138

$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws2.Range("F2").Value = 79

$ws4 = $wb.Worksheets.Item(4)   # 全部类型
$ws4.Range("F2").Value = 143
$ws4.Range("F3").Value = 79
$ws4.Range("F4").Value = 457
$ws4.Range("F7").Value = 11
$ws4.Range("F8").Value = 27
$ws4.Range("F10").Value = 138
